$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "46.124.10"
$ws.Range("E2").Value = "  -0.98%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.599.48"
$ws.Range("E3").Value = "  +0.83%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.09%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.64%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.98%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.596"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.65%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.09%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.581"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.33%  "

# Row 10 - Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.42%  "

# Row 11 - OKB
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.28"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.04%  "

# Row 12 - Dogecoin
$ws.Range("E12").Value = "  +0.21%  "

# Row 13 - Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.31%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.994.34"
$ws.Range("E14").Value = "  +0.70%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  +0.88%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.598.31"
$ws.Range("E16").Value = "  +0.49%  "

# Row 17 - Polygon
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.917"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.31%  "

# Row 18 - Chainlink
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.53%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "46.230.69"
$ws.Range("E19").Value = "  -1.07%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  +0.53%  "

# Row 21 - row21 -> Uniswap (swap)
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.18%  "

# Row 22 - row22 -> InternetComputer(DFINITY) (swap)
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.54%  "

# Row 23 - Litecoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.81%  "

# Row 24 - BitcoinCash
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "273.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.04%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +2.44%  "

# Row 26 - ImmutableX
$ws.Range("E26").Value = "  +0.55%  "

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.74%  "

# Row 28 - Dai
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.07%  "

# Row 29 - LEO
$ws.Range("E29").Value = "  +1.13%  "

# Row 30 - Cosmos
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.06%  "

# Row 31 - InjectiveProtocol
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.97%  "

# Row 32 - Toncoin
$ws.Range("E32").Value = "  -3.46%  "

# Row 33 - Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.67%  "

# Row 34 - LidoDAOToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.75%  "

# Row 35 - Monero
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "155.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.71%  "

# Row 36 - ARBITRUM
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.48%  "

# Row 37 - Hedera
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0837"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.01%  "

# Row 38 - WEMIXToken
$ws.Range("E38").Value = "  -4.24%  "

# Row 39 - Kaspa
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.124"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.39%  "

# Row 40 - Stellar
$ws.Range("E40").Value = "  +0.88%  "

# Row 41 - EnergySwap
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +26.26%  "

# Row 42 - Celestia
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.65%  "

# Row 43 - VeChain
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0331"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.43%  "

# Row 44 - NEARProtocol
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.10%  "

# Row 45 - RenderToken
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.26%  "

# Row 46 - Maker
$ws.Range("D46").Value = "2.100.98"
$ws.Range("E46").Value = "  +4.00%  "

# Row 47 - FirstDigitalUSD
$ws.Range("E47").Value = "  -0.05%  "

# Row 48 - BitcoinSV
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "95.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.78%  "

# Row 49 - FraxShare
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.62%  "

# Row 50 - Aave
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "108.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.07%  "

# Row 51 - row51 -> Algorand (replace Stacks)
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.200"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.05%  "
